$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.132470666666666
$ws.Range("H2").Value = 21.397412
$ws.Range("I2").Value = 0.1078130252899183
$ws.Range("J2").Value = 0.1078130252899183
$ws.Range("M2").Value = 0.9705896666666667
$ws.Range("N2").Value = 2.911769
$ws.Range("O2").Value = 0.02073452941466921
$ws.Range("P2").Value = 0.02073452941466921
$ws.Range("Q2").Value = 6.922702326869778
$ws.Range("R2").Value = 62.304320941828
$ws.Range("S2").Value = 0.002235452344158288
$ws.Range("T2").Value = 0.002235452344158288
$ws.Range("G3").Value = 7.132470666666666
$ws.Range("H3").Value = 21.397412
$ws.Range("I3").Value = 0.1078130252899183
$ws.Range("J3").Value = 0.1078130252899183
$ws.Range("O3").Value = 0.5628689972673966
$ws.Range("P3").Value = 0.5628689972673966
$ws.Range("Q3").Value = 187.92683639828
$ws.Range("R3").Value = 1691.34152758452
$ws.Range("S3").Value = 0.0606846094373008
$ws.Range("T3").Value = 0.0606846094373008
$ws.Range("G4").Value = 7.132470666666666
$ws.Range("H4").Value = 21.397412
$ws.Range("I4").Value = 0.1078130252899183
$ws.Range("J4").Value = 0.1078130252899183
$ws.Range("M4").Value = 19.49164633333333
$ws.Range("N4").Value = 58.47493899999999
$ws.Range("O4").Value = 0.4163964733179342
$ws.Range("P4").Value = 0.4163964733179341
$ws.Range("Q4").Value = 139.0235957175409
$ws.Range("R4").Value = 1251.212361457868
$ws.Range("S4").Value = 0.04489296350845924
$ws.Range("T4").Value = 0.04489296350845923
$ws.Range("I5").Value = 0.2490596131114117
$ws.Range("J5").Value = 0.2490596131114118
$ws.Range("M5").Value = 0.9705896666666667
$ws.Range("N5").Value = 2.911769
$ws.Range("O5").Value = 0.02073452941466921
$ws.Range("P5").Value = 0.02073452941466921
$ws.Range("Q5").Value = 15.99218237851345
$ws.Range("R5").Value = 143.929641406621
$ws.Range("S5").Value = 0.0051641338740647
$ws.Range("T5").Value = 0.005164133874064701
$ws.Range("I6").Value = 0.2490596131114117
$ws.Range("J6").Value = 0.2490596131114118
$ws.Range("O6").Value = 0.5628689972673966
$ws.Range("P6").Value = 0.5628689972673966
$ws.Range("S6").Value = 0.1401879346918261
$ws.Range("T6").Value = 0.1401879346918261
$ws.Range("I7").Value = 0.2490596131114117
$ws.Range("J7").Value = 0.2490596131114118
$ws.Range("M7").Value = 19.49164633333333
$ws.Range("N7").Value = 58.47493899999999
$ws.Range("O7").Value = 0.4163964733179342
$ws.Range("P7").Value = 0.4163964733179341
$ws.Range("Q7").Value = 321.1593670584613
$ws.Range("R7").Value = 2890.434303526151
$ws.Range("S7").Value = 0.103707544545521
$ws.Range("T7").Value = 0.103707544545521
$ws.Range("G8").Value = 42.546687
$ws.Range("H8").Value = 127.640061
$ws.Range("I8").Value = 0.6431273615986699
$ws.Range("J8").Value = 0.6431273615986699
$ws.Range("M8").Value = 0.9705896666666667
$ws.Range("N8").Value = 2.911769
$ws.Range("O8").Value = 0.02073452941466921
$ws.Range("P8").Value = 0.02073452941466921
$ws.Range("Q8").Value = 41.295374753101
$ws.Range("R8").Value = 371.658372777909
$ws.Range("S8").Value = 0.01333494319644622
$ws.Range("T8").Value = 0.01333494319644622
$ws.Range("G9").Value = 42.546687
$ws.Range("H9").Value = 127.640061
$ws.Range("I9").Value = 0.6431273615986699
$ws.Range("J9").Value = 0.6431273615986699
$ws.Range("O9").Value = 0.5628689972673966
$ws.Range("P9").Value = 0.5628689972673966
$ws.Range("Q9").Value = 1121.02308734409
$ws.Range("R9").Value = 10089.20778609681
$ws.Range("S9").Value = 0.3619964531382697
$ws.Range("T9").Value = 0.3619964531382697
$ws.Range("G10").Value = 42.546687
$ws.Range("H10").Value = 127.640061
$ws.Range("I10").Value = 0.6431273615986699
$ws.Range("J10").Value = 0.6431273615986699
$ws.Range("M10").Value = 19.49164633333333
$ws.Range("N10").Value = 58.47493899999999
$ws.Range("O10").Value = 0.4163964733179342
$ws.Range("P10").Value = 0.4163964733179341
$ws.Range("Q10").Value = 829.3049756590309
$ws.Range("R10").Value = 7463.744780931277
$ws.Range("S10").Value = 0.2677959652639539
$ws.Range("T10").Value = 0.2677959652639539
